$d = $word.ActiveDocument
$d.Content.Find.Execute("????", $false, $false, $false, $false, $false, $true, 1, $false, "30th Sep 2022", 2)

# Now clear highlight on the replaced text
$r = $d.Content
$r.Find.Execute("30th Sep 2022", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.HighlightColorIndex = 0
